$d = $word.ActiveDocument

# Before state: a single empty paragraph that only contains the leftover
# "_GoBack" bookmark pair (Word drops one of these at the last edit point).
#
# Target state (per the diff):
#   Para 1: "Meus estudos de GIT"
#   Para 2: "Meu primeiro exercício com GIT Comm" + <bookmark _GoBack/> + "and"
#           (i.e. the visible text reads "...GIT Command", but the original
#           "_GoBack" bookmark stays sandwiched between the two runs, exactly
#           where the author's cursor was while typing.)

# 1) Add the new first paragraph in front of everything else.
$start = $d.Range(0, 0)
$start.InsertBefore("Meus estudos de GIT`r")

# 2) The "_GoBack" bookmark now lives in the second paragraph. Type "and"
#    right at the bookmark's position first: Word places freshly typed text
#    *after* the bookmark markers, which is what keeps bookmarkStart/End
#    together and ahead of this run.
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks("_GoBack")
    $sel = $word.Selection
    $sel.SetRange($bm.Start, $bm.Start)
    $sel.TypeText("and")

    # 3) Now insert the lead-in text directly before the (still zero-width)
    #    bookmark location, using Range.InsertBefore so it lands ahead of
    #    the bookmarkStart element instead of merging into the "and" run.
    $bm2 = $d.Bookmarks("_GoBack")
    $rngBefore = $d.Range($bm2.Start, $bm2.Start)
    $rngBefore.InsertBefore("Meu primeiro exercício com GIT Comm")
} else {
    # Fallback (shouldn't happen for this document): just write the full
    # sentence into the second paragraph.
    $p2 = $d.Paragraphs(2).Range
    $p2.InsertBefore("Meu primeiro exercício com GIT Command")
}
